# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types) sheets,
# matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 816
$ws1.Range("F4").Value = 305
$ws1.Range("F5").Value = 512
$ws1.Range("F6").Value = 1149
$ws1.Range("F11").Value = 1192
$ws1.Range("F14").Value = 871
$ws1.Range("F15").Value = 865
$ws1.Range("F17").Value = 68
$ws1.Range("F20").Value = 761
$ws1.Range("F21").Value = 1737
$ws1.Range("F22").Value = 2798
$ws1.Range("F23").Value = 808
$ws1.Range("F25").Value = 2147
$ws1.Range("F26").Value = 666
$ws1.Range("F27").Value = 2988
$ws1.Range("F28").Value = 576
$ws1.Range("F31").Value = 88
$ws1.Range("F32").Value = 724
$ws1.Range("F33").Value = 145
$ws1.Range("F34").Value = 129
$ws1.Range("F36").Value = 1054
$ws1.Range("F37").Value = 1758
$ws1.Range("F38").Value = 380
$ws1.Range("F41").Value = 182
$ws1.Range("F44").Value = 41

# --- 演出 (sheet 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 76

# --- 全部类型 (sheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 816
$ws4.Range("F4").Value = 305
$ws4.Range("F5").Value = 512
$ws4.Range("F6").Value = 1149
$ws4.Range("F11").Value = 1192
$ws4.Range("F13").Value = 871
$ws4.Range("F14").Value = 865
$ws4.Range("F17").Value = 68
$ws4.Range("F20").Value = 761
$ws4.Range("F21").Value = 1737
$ws4.Range("F22").Value = 2798
$ws4.Range("F23").Value = 808
$ws4.Range("F26").Value = 2988
$ws4.Range("F27").Value = 576
$ws4.Range("F34").Value = 88
$ws4.Range("F35").Value = 76
$ws4.Range("F36").Value = 724
$ws4.Range("F37").Value = 145
$ws4.Range("F38").Value = 129
$ws4.Range("F41").Value = 1054
$ws4.Range("F42").Value = 1758
$ws4.Range("F43").Value = 380
$ws4.Range("F45").Value = 182
$ws4.Range("F48").Value = 41
